$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that held the per-fuel breakdown (hydrogen/biomass rows);
# the data for all fuels now lives in a single row (row 2).
$ws.Rows("3:4").Delete()

# Extend the header row/style formatting from B1 out to L1, reusing the
# existing bold/bordered/centered header style instead of creating a new one.
$ws.Range("B1").Copy()
$ws.Range("E1:L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update header labels (B1:D1 relabeled, E1:L1 newly populated).
$ws.Range("B1").Value = "Total Cost"
$ws.Range("C1").Value = "crudeoil"
$ws.Range("D1").Value = "hydrogen"
$ws.Range("E1").Value = "biomass"
$ws.Range("F1").Value = "RefineryProduction"
$ws.Range("G1").Value = "MtGProduction"
$ws.Range("H1").Value = "GtkmProduction"
$ws.Range("I1").Value = "B2gasProduction"
$ws.Range("J1").Value = "GasHubUsage"
$ws.Range("K1").Value = "KmHubUsage"
$ws.Range("L1").Value = "KilometersUsage"

# Update the data row. A2 (index column) is unchanged.
$ws.Range("B2").Value = 124.024199843872
$ws.Range("C2").Value = 546.448087431694
$ws.Range("D2").Value = 2371.1943793911
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 508.1967213114755
$ws.Range("G2").Value = 1991.803278688524
$ws.Range("H2").Value = 1000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 2500
$ws.Range("K2").Value = 1000
$ws.Range("L2").Value = 1000
